$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 508, shifting all existing rows
# 508..624 down to 510..626 (keeps their values/formatting intact).
$ws.Rows.Item(508).Resize(2).Insert()

# --- New row 508 ---
$ws.Range("A508").Value = 1
$ws.Range("B508").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C508").Value = "Arica y Parinacota"
$ws.Range("D508").Value = 45173
$ws.Range("E508").Value = 15
$ws.Range("F508").Value = 100112023
$ws.Range("G508").Value = "Brócoli"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Segunda"
$ws.Range("J508").Value = 1200
$ws.Range("K508").Value = 400
$ws.Range("L508").Value = 500
$ws.Range("M508").Value = 450
$ws.Range("N508").Value = "$/unidad"
$ws.Range("O508").Value = "Región de Arica y Parinacota"
$ws.Range("P508").Value = 450
$ws.Range("Q508").Value = 1
$ws.Range("R508").Value = "Hortaliza"

# --- New row 509 ---
$ws.Range("A509").Value = 1
$ws.Range("B509").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C509").Value = "Arica y Parinacota"
$ws.Range("D509").Value = 45173
$ws.Range("E509").Value = 15
$ws.Range("F509").Value = 100112023
$ws.Range("G509").Value = "Brócoli"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Tercera"
$ws.Range("J509").Value = 1200
$ws.Range("K509").Value = 200
$ws.Range("L509").Value = 300
$ws.Range("M509").Value = 250
$ws.Range("N509").Value = "$/unidad"
$ws.Range("O509").Value = "Región de Arica y Parinacota"
$ws.Range("P509").Value = 250
$ws.Range("Q509").Value = 1
$ws.Range("R509").Value = "Hortaliza"
